$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Column B (First name initials) for every student row, cycling
# through a, b, c, d. The first four assignments are ordered a, d,
# b, c so that the workbook's shared string table picks up the new
# strings in that same order (indices 10=a, 11=d, 12=b, 13=c), just
# like the target workbook.
# ------------------------------------------------------------------
$ws.Range("B2").Value = "a"
$ws.Range("B5").Value = "d"
$ws.Range("B3").Value = "b"
$ws.Range("B4").Value = "c"
$ws.Range("B6").Value = "a"
$ws.Range("B7").Value = "b"
$ws.Range("B8").Value = "c"
$ws.Range("B9").Value = "d"
$ws.Range("B10").Value = "a"
$ws.Range("B11").Value = "b"
$ws.Range("B12").Value = "c"
$ws.Range("B13").Value = "d"
$ws.Range("B14").Value = "a"
$ws.Range("B15").Value = "b"
$ws.Range("B16").Value = "c"
$ws.Range("B17").Value = "d"
$ws.Range("B18").Value = "a"
$ws.Range("B19").Value = "b"
$ws.Range("B20").Value = "c"

# ------------------------------------------------------------------
# Column D (Eng) and Column E (Maths) marks for every student row.
# Column G (Total) is computed automatically via the existing
# SUM(D:F) formulas already present in the sheet.
# ------------------------------------------------------------------
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 45

$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 32

$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 45

$ws.Range("D5").Value = 45
$ws.Range("E5").Value = 45

$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 7

$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 9

$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 36

$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 14

$ws.Range("D10").Value = 48
$ws.Range("E10").Value = 48

$ws.Range("D11").Value = 42
$ws.Range("E11").Value = 42

$ws.Range("D12").Value = 36
$ws.Range("E12").Value = 36

$ws.Range("D13").Value = 38
$ws.Range("E13").Value = 38

$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 7

$ws.Range("D15").Value = 9
$ws.Range("E15").Value = 9

$ws.Range("D16").Value = 36
$ws.Range("E16").Value = 36

$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 14

$ws.Range("D18").Value = 48
$ws.Range("E18").Value = 48

$ws.Range("D19").Value = 42
$ws.Range("E19").Value = 42

$ws.Range("D20").Value = 36
$ws.Range("E20").Value = 36

# ------------------------------------------------------------------
# Update the active cell / selection to match the author's final
# cursor position in the sheet.
# ------------------------------------------------------------------
$ws.Range("L8").Select()
